# "Generate Report for Handoff"
#
# The localization pipeline re-ran for the "b.md" source file: a new handoff
# package was generated (status flips from "Handed back: in sync with en-US"
# to "Ready for handoff"), a new xliff handoff file + timestamp is recorded,
# and because the previously handed-back file is now stale relative to the
# new handoff, an Error Detail note is attached. This touches the Overview
# sheet's "b.md" row plus the per-locale (zh-cn / de-de) detail sheets' "b.md"
# row.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe0ec556dc0675ab517216d6da5d32e898bf5cf2/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a76d2de706d57dc5af0315d60eb2978cbdd525/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" file. Status columns (zh-cn / de-de)
# move to "Ready for handoff" and the "Latest HO Xliff Generate Date" bumps
# to the new handoff timestamp.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-21 10:44:31"

# ---------------------------------------------------------------------
# zh-cn detail sheet: row 3 ("b.md") gets the new handoff file/timestamp,
# flips Status to "Ready for handoff", and records the stale-handback
# Error Detail. Column P (Error Detail) widens to fit the new text.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-21 10:44:27"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de detail sheet: same shape of update as zh-cn, using the de-de
# handoff file name / the Overview-sheet timestamp.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-21 10:44:31"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
